$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Module 3 (Privacy/Ethics) section: rename "Belmont principles" to
# "Ethical principles", reorder terms, and tweak the beneficence definition ---

# Row 18 (Confidentiality) / Row 19 (Privacy): swap their definitions
$ws.Range("C18").Value = "The duty of anyone entrusted with health information to keep that information private"
$ws.Range("C19").Value = "The right of an individual to keep his or her information (health related or otherwise) private"

# Row 20: Ethical principles: respect for persons
$ws.Range("B20").Value = "Ethical principles: respect for persons"
$ws.Range("C20").Value = "Defined by two ethical convictions: a) individuals should be treated as autonomous agents; b) persons with diminished autonomy are entitled to protection"

# Row 21: Ethical principles: justice
$ws.Range("B21").Value = "Ethical principles: justice"
$ws.Range("C21").Value = "Ethical principle that the burdens and benefits of research and public health practice should be justly distributed, including attention to need, effort, contribution, and merit"

# Row 22: Ethical principles: beneficence
$ws.Range("B22").Value = "Ethical principles: beneficence"
$ws.Range("C22").Value = "Two general rules have been formulated as complementary expressions of beneficent actions in this sense: (1) do not harm (e.g. non-maleficence) and (2) maximize possible benefits and minimize possible harms"

# Row 23: Geomask
$ws.Range("B23").Value = "Geomask"
$ws.Range("C23").Value = "A class of methods for changing the geographic location of an individual in an unpredictable way to protect confidentiality, while trying to preserve the relationship between geocoded locations and disease occurrence (Sherman and Fetters 2007, Wiggins 2002)"

# --- New Module 4: "Disease Mapping 1" content (rows 24-30) ---

$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Spatial heterogeneity"
$ws.Range("C24").Value = "Attributes or statistical parameters are varied (e.g. not homogenous) across sub-areas in a broader region. In Disease mapping we typically are evaluating whether (and how much) disease intensity (risk, rate, prevalence) varies across places. "

$ws.Range("A25").Value = 4
$ws.Range("B25").Value = "Spatial dependence"
$ws.Range("C25").Value = "When attribute values or statistical parameters are, on avreage, more similar for nearby places than they are for distant places. Spatial dependence is evaluated by looking at pairs or sets of places."

$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Stationarity vs non-stationarity"
$ws.Range("C26").Value = "Many statistics assume that the parameter, estimate, or property is constant across sub-units. For example if we take the average height of a population, under stationarity we would assume that average applies equally to sub-populations. In contrast, non-stationarity implies the parameter, estimate, or property varies across sub-groups. In spatial analysis stationarity is an assumption of homogeneity, and non-stationarity allows for heterogeneity."

$ws.Range("A27").Value = 4
$ws.Range("B27").Value = "Global vs Local spatial analysis"
$ws.Range("C27").Value = "Global analysis evaluates a pattern or trends that characterizes the entire study region; in contrast local analysis characterizes patterns that are unique to each sub-region of the study area"

$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Standardize Morbidity/Mortality Ratio (SMR)"
$ws.Range("C28").Value = "The ratio of observed to expected disease morbidity or mortality. Often the 'expected' is defined as the overall population (or study-specific) rate; in that case the SMR indicates the relative deviation of a specific unit from the global or overall rate"

$ws.Range("A29").Value = 4
$ws.Range("B29").Value = "Bayesian methods"
$ws.Range("C29").Value = "Methods of statistical inference in which Bayes' theorem is used to update the probability for a hypothesis as more evidence or information becomes available. In disease mapping, the Bayesian framework is frequently used to accomplish rate stabilization and smoothing by using global or local data to inform the 'prior' "

$ws.Range("A30").Value = 4
$ws.Range("B30").Value = "Empirical Bayes methods"
$ws.Range("C30").Value = "Estimation procedures in a Bayesian framework in which the prior distribution is estimated from the data. In disease mapping, Empirical Bayes estimators use global or local disease information as a prior in estimating (and smoothing/stabilizing) each local rate"

# Update the selected/active cell to match the edited workbook
$ws.Range("C27").Select()
